{"js": "// Replace the date heading and all twenty-five \"A\u00d7B=C\" answer cells in the\n// practice-sheet table with the updated values from the target revision.\n// Each entry is [oldText, newText]; every old value is a unique substring\n// of the document body, so a plain search + insertText(replace) on each\n// hit swaps the text while preserving that run's existing formatting\n// (font, size, etc.).\nconst replacements = [\n  [\"2025-01-20 Monday\", \"2025-01-21 Tuesday\"],\n  [\"191\u00d75=955\", \"507\u00d75=2535\"],\n  [\"976\u00d72=1952\", \"674\u00d77=4718\"],\n  [\"137\u00d75=685\", \"763\u00d79=6867\"],\n  [\"872\u00d75=4360\", \"444\u00d77=3108\"],\n  [\"296\u00d74=1184\", \"569\u00d74=2276\"],\n  [\"612\u00d73=1836\", \"704\u00d77=4928\"],\n  [\"857\u00d74=3428\", \"603\u00d79=5427\"],\n  [\"455\u00d79=4095\", \"542\u00d76=3252\"],\n  [\"820\u00d73=2460\", \"594\u00d72=1188\"],\n  [\"670\u00d75=3350\", \"146\u00d78=1168\"],\n  [\"559\u00d73=1677\", \"933\u00d78=7464\"],\n  [\"473\u00d76=2838\", \"149\u00d74=596\"],\n  [\"192\u00d74=768\", \"291\u00d79=2619\"],\n  [\"314\u00d76=1884\", \"402\u00d75=2010\"],\n  [\"319\u00d77=2233\", \"127\u00d79=1143\"],\n  [\"508\u00d77=3556\", \"246\u00d72=492\"],\n  [\"888\u00d79=7992\", \"679\u00d75=3395\"],\n  [\"961\u00d78=7688\", \"822\u00d73=2466\"],\n  [\"554\u00d77=3878\", \"145\u00d77=1015\"],\n  [\"765\u00d72=1530\", \"487\u00d72=974\"],\n  [\"783\u00d76=4698\", \"996\u00d73=2988\"],\n  [\"530\u00d76=3180\", \"108\u00d77=756\"],\n  [\"117\u00d73=351\", \"786\u00d78=6288\"],\n  [\"819\u00d78=6552\", \"551\u00d74=2204\"],\n  [\"738\u00d79=6642\", \"966\u00d76=5796\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the date heading and all twenty-five \"A x B = C\" answer cells in\n# the practice-sheet table with the updated values from the target revision.\n# Each entry is (oldText, newText); every old value is a unique substring of\n# the document, so Find/Replace on each pair swaps the text in place while\n# leaving the owning run's formatting (font, size, etc.) untouched.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2025-01-20 Monday\", \"2025-01-21 Tuesday\"),\n    @(\"191\u00d75=955\", \"507\u00d75=2535\"),\n    @(\"976\u00d72=1952\", \"674\u00d77=4718\"),\n    @(\"137\u00d75=685\", \"763\u00d79=6867\"),\n    @(\"872\u00d75=4360\", \"444\u00d77=3108\"),\n    @(\"296\u00d74=1184\", \"569\u00d74=2276\"),\n    @(\"612\u00d73=1836\", \"704\u00d77=4928\"),\n    @(\"857\u00d74=3428\", \"603\u00d79=5427\"),\n    @(\"455\u00d79=4095\", \"542\u00d76=3252\"),\n    @(\"820\u00d73=2460\", \"594\u00d72=1188\"),\n    @(\"670\u00d75=3350\", \"146\u00d78=1168\"),\n    @(\"559\u00d73=1677\", \"933\u00d78=7464\"),\n    @(\"473\u00d76=2838\", \"149\u00d74=596\"),\n    @(\"192\u00d74=768\", \"291\u00d79=2619\"),\n    @(\"314\u00d76=1884\", \"402\u00d75=2010\"),\n    @(\"319\u00d77=2233\", \"127\u00d79=1143\"),\n    @(\"508\u00d77=3556\", \"246\u00d72=492\"),\n    @(\"888\u00d79=7992\", \"679\u00d75=3395\"),\n    @(\"961\u00d78=7688\", \"822\u00d73=2466\"),\n    @(\"554\u00d77=3878\", \"145\u00d77=1015\"),\n    @(\"765\u00d72=1530\", \"487\u00d72=974\"),\n    @(\"783\u00d76=4698\", \"996\u00d73=2988\"),\n    @(\"530\u00d76=3180\", \"108\u00d77=756\"),\n    @(\"117\u00d73=351\", \"786\u00d78=6288\"),\n    @(\"819\u00d78=6552\", \"551\u00d74=2204\"),\n    @(\"738\u00d79=6642\", \"966\u00d76=5796\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n\n    $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n"}
